# ---------------------------------------------------------------------------
# Scheduled runner update: refresh market-price-derived columns (H:N) on the
# Leve profit tables for each job sheet (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
#
# Columns:
#   H = currentAveragePrice     I = currentAveragePriceNQ
#   J = currentAveragePriceHQ   K = LevePriceNQ
#   L = LevePriceHQ             M = LeveProfitNQ   N = LeveProfitHQ
#
# Values are plain refreshed data points (not formulas), written directly.
# An empty string clears a cell entirely (removes it from the sheet), used
# where a column no longer applies for that leve row.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")

# Row 15
$ws.Range("H15").Value = 1981.5
$ws.Range("I15").Value = 1981.5
$ws.Range("K15").Value = 5944.5
$ws.Range("M15").Value = -5775.5

# Row 16
$ws.Range("H16").Value = 3858.5715
$ws.Range("J16").Value = 4202
$ws.Range("L16").Value = 4202
$ws.Range("N16").Value = -4662

# Row 69
$ws.Range("H69").Value = 18084.334
$ws.Range("I69").Value = 7664
$ws.Range("J69").Value = 28504.666
$ws.Range("K69").Value = 22992
$ws.Range("L69").Value = 85513.99800000001
$ws.Range("M69").Value = -22118
$ws.Range("N69").Value = -87261.99800000001

# Row 72
$ws.Range("H72").Value = 18084.334
$ws.Range("I72").Value = 7664
$ws.Range("J72").Value = 28504.666
$ws.Range("K72").Value = 68976
$ws.Range("L72").Value = 256541.994
$ws.Range("M72").Value = -64608
$ws.Range("N72").Value = -265277.994

# Row 107
$ws.Range("H107").Value = 1406.8
$ws.Range("I107").Value = 1508.5
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 1508.5
$ws.Range("L107").Value = 1000
$ws.Range("M107").Value = 411.5
$ws.Range("N107").Value = -4840

# Row 132
$ws.Range("H132").Value = 2559.8838
$ws.Range("I132").Value = 2141.2896
$ws.Range("K132").Value = 6423.8688
$ws.Range("M132").Value = -3893.8688

# Row 141
$ws.Range("H141").Value = 2412.1428
$ws.Range("I141").Value = 2412.1428
$ws.Range("K141").Value = 7236.428400000001
$ws.Range("M141").Value = -2056.428400000001


# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")

# Row 31
$ws.Range("H31").Value = 34427
$ws.Range("J31").Value = 111960
$ws.Range("L31").Value = 111960
$ws.Range("N31").Value = -112548

# Row 92
$ws.Range("H92").Value = 32497.5
$ws.Range("J92").Value = 19996.666
$ws.Range("L92").Value = 19996.666
$ws.Range("N92").Value = -24988.666

# Row 122
$ws.Range("H122").Value = 1200
$ws.Range("I122").Value = 1200
$ws.Range("K122").Value = 3600
$ws.Range("M122").Value = -1150

# Row 123
$ws.Range("H123").Value = 88000
$ws.Range("J123").Value = 88000
$ws.Range("L123").Value = 88000
$ws.Range("N123").Value = -97800

# Row 125
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = ""
$ws.Range("N125").Value = 0


# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")

# Row 20
$ws.Range("H20").Value = 3918.5334
$ws.Range("I20").Value = 3798.6365
$ws.Range("K20").Value = 3798.6365
$ws.Range("M20").Value = -3551.6365

# Row 28
$ws.Range("H28").Value = 46970
$ws.Range("J28").Value = 46970
$ws.Range("L28").Value = 46970
$ws.Range("N28").Value = -47558

# Row 96
$ws.Range("H96").Value = 39262.453
$ws.Range("J96").Value = 71097.8
$ws.Range("L96").Value = 71097.8
$ws.Range("N96").Value = -76589.8


# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")

# Row 51
$ws.Range("H51").Value = 54999.832
$ws.Range("I51").Value = 25000
$ws.Range("K51").Value = 25000
$ws.Range("M51").Value = -24264

# Row 58
$ws.Range("H58").Value = 1000
$ws.Range("I58").Value = 1000
$ws.Range("K58").Value = 1000
$ws.Range("M58").Value = -797

# Row 61
$ws.Range("H61").Value = 54999.832
$ws.Range("I61").Value = 25000
$ws.Range("K61").Value = 25000
$ws.Range("M61").Value = -24652

# Row 62
$ws.Range("H62").Value = 3992.4285
$ws.Range("I62").Value = 3993.4
$ws.Range("J62").Value = 3990
$ws.Range("K62").Value = 3993.4
$ws.Range("L62").Value = 3990
$ws.Range("M62").Value = -3369.4
$ws.Range("N62").Value = -5238

# Row 65
$ws.Range("H65").Value = 3992.4285
$ws.Range("I65").Value = 3993.4
$ws.Range("J65").Value = 3990
$ws.Range("K65").Value = 19967
$ws.Range("L65").Value = 19950
$ws.Range("M65").Value = -16847
$ws.Range("N65").Value = -26190

# Row 103
$ws.Range("H103").Value = 34972.25
$ws.Range("I103").Value = 24999.5
$ws.Range("J103").Value = 44945
$ws.Range("K103").Value = 24999.5
$ws.Range("L103").Value = 44945
$ws.Range("M103").Value = -23827.5
$ws.Range("N103").Value = -47289

# Row 132
$ws.Range("H132").Value = 3641.6365
$ws.Range("I132").Value = 3641.6365
$ws.Range("K132").Value = 10924.9095
$ws.Range("M132").Value = -8394.9095

# Row 136
$ws.Range("H136").Value = 1000
$ws.Range("I136").Value = 1000
$ws.Range("K136").Value = 3000
$ws.Range("M136").Value = -450


# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")

# Row 16
$ws.Range("H16").Value = 1954.3334
$ws.Range("I16").Value = 298.14285
$ws.Range("J16").Value = 7751
$ws.Range("K16").Value = 894.4285500000001
$ws.Range("L16").Value = 23253
$ws.Range("M16").Value = -721.4285500000001
$ws.Range("N16").Value = -23599

# Row 25
$ws.Range("H25").Value = 1
$ws.Range("I25").Value = 1
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 3
$ws.Range("L25").Value = ""
$ws.Range("N25").Value = 0
$ws.Range("M25").Value = 166

# Row 30
$ws.Range("H30").Value = 1
$ws.Range("I30").Value = 1
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 3
$ws.Range("L30").Value = ""
$ws.Range("N30").Value = 0
$ws.Range("M30").Value = 99

# Row 44
$ws.Range("H44").Value = 18128.084
$ws.Range("I44").Value = 21781.889
$ws.Range("J44").Value = 7166.6665
$ws.Range("K44").Value = 65345.667
$ws.Range("L44").Value = 21499.9995
$ws.Range("M44").Value = -64947.667
$ws.Range("N44").Value = -22295.9995

# Row 131
$ws.Range("H131").Value = 6350.3687
$ws.Range("I131").Value = 7268.6665
$ws.Range("J131").Value = 5523.9
$ws.Range("K131").Value = 21805.9995
$ws.Range("L131").Value = 16571.7
$ws.Range("M131").Value = -16765.9995
$ws.Range("N131").Value = -26651.7


# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")

# Row 102
$ws.Range("H102").Value = 9437.375
$ws.Range("I102").Value = 7333.222
$ws.Range("J102").Value = 12142.714
$ws.Range("K102").Value = 7333.222
$ws.Range("L102").Value = 12142.714
$ws.Range("M102").Value = -5711.222
$ws.Range("N102").Value = -15386.714

# Row 122
$ws.Range("H122").Value = 1942.4286
$ws.Range("I122").Value = 1500
$ws.Range("J122").Value = 2016.1666
$ws.Range("K122").Value = 4500
$ws.Range("L122").Value = 6048.4998
$ws.Range("M122").Value = -2050
$ws.Range("N122").Value = -10948.4998

# Row 136
$ws.Range("H136").Value = 38647.75
$ws.Range("J136").Value = 38647.75
$ws.Range("L136").Value = 115943.25
$ws.Range("N136").Value = -121043.25


# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")

# Row 22
$ws.Range("H22").Value = 2852.2173
$ws.Range("I22").Value = 2787.625
$ws.Range("K22").Value = 2787.625
$ws.Range("M22").Value = -2492.625

# Row 26
$ws.Range("H26").Value = 9715.294
$ws.Range("I26").Value = 5000
$ws.Range("J26").Value = 10010
$ws.Range("K26").Value = 5000
$ws.Range("L26").Value = 10010
$ws.Range("M26").Value = -4705
$ws.Range("N26").Value = -10600

# Row 27
$ws.Range("H27").Value = 2852.2173
$ws.Range("I27").Value = 2787.625
$ws.Range("K27").Value = 2787.625
$ws.Range("M27").Value = -2680.625

# Row 30
$ws.Range("H30").Value = 3607.5
$ws.Range("I30").Value = 3607.5
$ws.Range("K30").Value = 3607.5
$ws.Range("M30").Value = -3499.5

# Row 31
$ws.Range("H31").Value = 737.7143
$ws.Range("I31").Value = 637.5
$ws.Range("J31").Value = 871.3333
$ws.Range("K31").Value = 637.5
$ws.Range("L31").Value = 871.3333
$ws.Range("M31").Value = -389.5
$ws.Range("N31").Value = -1367.3333

# Row 33
$ws.Range("H33").Value = 3729.25
$ws.Range("I33").Value = 2450
$ws.Range("K33").Value = 2450
$ws.Range("M33").Value = -2160

# Row 35
$ws.Range("H35").Value = 1712
$ws.Range("I35").Value = 1712
$ws.Range("K35").Value = 1712
$ws.Range("M35").Value = -1376

# Row 46
$ws.Range("H46").Value = 4456
$ws.Range("I46").Value = 4118.8335
$ws.Range("K46").Value = 4118.8335
$ws.Range("M46").Value = -3930.8335

# Row 55
$ws.Range("H55").Value = 66666970
$ws.Range("I55").Value = 90909336
$ws.Range("K55").Value = 90909336
$ws.Range("M55").Value = -90909163

# Row 95
$ws.Range("H95").Value = 27249.5
$ws.Range("J95").Value = 27249.5
$ws.Range("L95").Value = 27249.5
$ws.Range("N95").Value = -32741.5

# Row 98
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = ""
$ws.Range("N98").Value = 0

# Row 102
$ws.Range("H102").Value = 79332.664
$ws.Range("J102").Value = 79000
$ws.Range("L102").Value = 79000
$ws.Range("N102").Value = -85490

# Row 122
$ws.Range("H122").Value = 5609.3477
$ws.Range("I122").Value = 4677.4116
$ws.Range("K122").Value = 14032.2348
$ws.Range("M122").Value = -11582.2348


# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")

# Row 58
$ws.Range("H58").Value = 43997.5
$ws.Range("J58").Value = 43997.5
$ws.Range("L58").Value = 43997.5
$ws.Range("N58").Value = -44613.5

# Row 118
$ws.Range("H118").Value = 79995
$ws.Range("J118").Value = 79995
$ws.Range("L118").Value = 79995
$ws.Range("N118").Value = -83309

# Row 132
$ws.Range("H132").Value = 10814.63
$ws.Range("I132").Value = 3149.5
$ws.Range("J132").Value = 16946.732
$ws.Range("K132").Value = 9448.5
$ws.Range("L132").Value = 50840.196
$ws.Range("M132").Value = -6918.5
$ws.Range("N132").Value = -55900.196
